# Insert a new salesperson row ("SARMIENTO SARMIENTO SANDRA EULALIA") at row 19
# on both worksheets, pushing the existing row 19 (VIEJO RIVAS MAYRA ANABELLE)
# and the totals row below it down by one, and refresh the "0 de N" counter text
# on the sheet that uses it.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Duplicate row 19 (copy + insert) so the new row inherits the same column
    # formatting/styles as the existing data rows, then push the old row 19 and
    # everything below it down by one row.
    $ws.Rows.Item(19).Copy()
    $ws.Rows.Item(19).Insert()

    # A19 already reads "ILLER LOPEZ ROBERTO FERNANDO" (copied from the row
    # above); only the client/customer name in column B needs to change.
    $ws.Cells.Item(19, 2).Value = "SARMIENTO SARMIENTO SANDRA EULALIA"

    # Determine how many data columns this sheet has by scanning the header row.
    $lastCol = 1
    while ($ws.Cells.Item(1, $lastCol + 1).Value2 -ne $null -and $ws.Cells.Item(1, $lastCol + 1).Value2 -ne "") {
        $lastCol = $lastCol + 1
    }

    # The totals row (now shifted to row 21) may use a "0 de N" textual counter;
    # bump it from 18 to 19 to account for the newly added row.
    $totalsText = $ws.Cells.Item(21, 3).Value2
    if ($totalsText -eq "0 de 18") {
        for ($c = 3; $c -le $lastCol; $c++) {
            $ws.Cells.Item(21, $c).Value = "0 de 19"
        }
    }
}
